# Update cryptocurrency Price and Volume(1h) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.987.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.021.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.31%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "547.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.02%  "

$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.013.42"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.44%  "

$ws.Range("E9").Value = "  +0.26%  "

$ws.Range("E10").Value = "  -4.43%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.14"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.66%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.451"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.44%  "

$ws.Range("E13").Value = "  -1.86%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.513.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.27%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.065.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.11%  "

$ws.Range("E17").Value = "  -2.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.024.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.14%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "478.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.28%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.87%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.677"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.61%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.51%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("E27").Value = "  +0.78%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.29%  "

$ws.Range("E29").Value = "  +0.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.88%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.15%  "

$ws.Range("E32").Value = "  +1.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.32"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.91%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "55.47"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.94%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.95"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.92%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "459.56"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.87%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.234.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0799"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.38%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0386"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.23%  "

$ws.Range("E41").Value = "  +0.37%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.43%  "

$ws.Range("E43").Value = "  -7.37%  "

$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("E45").Value = "  -3.71%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.32%  "

$ws.Range("E47").Value = "  -2.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.108"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.42%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "118.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.55%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₃0497"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.48%  "

$ws.Range("E51").Value = "  +6.97%  "
